$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 107; $row++) {
    $cell = $ws.Cells.Item($row, 4)
    $val = $cell.Value()
    if ($val -eq "F") {
        $cell.Value = "Female"
    } elseif ($val -eq "M") {
        $cell.Value = "Male"
    }
}
